# Gestor de actuaciones - olin-3Incidencias.xlsx
# "Ajusto los css de la muestra de actividades pendientes de coordinar"
#
# Adds the "Poblacion Instalacion" (Q) / "Direccion Instalacion" (R) data
# for the three pending-coordination rows, and leaves the selection on the
# next empty cell (S4) the way it was left after typing the last address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Pepito Perez / Malaga
$ws.Range("Q2").Value2 = "Málaga"
$ws.Range("R2").Value2 = "C/ Málaga, 45, 5ºB"

# Row 3 - Manolo el del tambor / Fuengirola
$ws.Range("Q3").Value2 = "Fuengirola"
$ws.Range("R3").Value2 = "C/ Fuengirola, Urbanización el coto. 46"

# Row 4 - Fernando Fernandez / Mijas costa
$ws.Range("Q4").Value2 = "Mijas costa"
$ws.Range("R4").Value2 = "Urbanización el Candado, calle Fuensanta, 45"

# Scroll the view so the newly filled columns are visible and leave the
# selection where the user ended up after filling in the last address.
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$selected = $ws.Range("S4").Select()
